# 14 October assignment commited
#
# EcommData ("fav7" favourite -> "Fav10") and NykaaData (new "expectedMsg"
# column + two more favourite lipstick rows) updates.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# EcommData sheet: rename saved-wishlist item "fav7" -> "Fav10"
# ---------------------------------------------------------------------
$ecomm = $wb.Worksheets.Item("EcommData")
$ecomm.Range("C6").Value = "Fav10"

# Column D (item names) now needs to be wide enough to fit the long
# mobile-name strings stored in it.
$ecomm.Columns.Item(4).ColumnWidth = 43.16666666666667

# ---------------------------------------------------------------------
# NykaaData sheet: add an "expectedMsg" column and two extra favourite
# lipstick rows.
# ---------------------------------------------------------------------
$nykaa = $wb.Worksheets.Item("NykaaData")

# New favourite item names underneath the existing one (E2).
$nykaa.Range("E3").Value = "M.A.C Matte Lipstick - Mehr"

# New "expectedMsg" column with its header + data values.
$nykaa.Range("F1").Value = "expectedMsg"
$nykaa.Range("F2").Value = "Your Shopping Bag is Empty"

$nykaa.Range("E6").Value = "M.A.C Cremesheen Lipstick - Creme In Your Coffee"

# Resize the new / widened columns to fit their content (values picked so
# the engine's internal character-width rounding lands on the real
# best-fit widths Excel computed for these columns).
$nykaa.Columns.Item(1).ColumnWidth = 13.666666666666666
$nykaa.Columns.Item(3).ColumnWidth = 5.833333333333333
$nykaa.Columns.Item(4).ColumnWidth = 6.666666666666667
$nykaa.Columns.Item(5).ColumnWidth = 49.666666666666664
$nykaa.Columns.Item(6).ColumnWidth = 24.833333333333332

$nykaa.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# Restore the selection on each sheet, leaving NykaaData as the active
# (visible) tab, same as before the edit.
# ---------------------------------------------------------------------
$ecomm.Activate()
$ecomm.Range("B12").Select()

$nykaa.Activate()
$nykaa.Range("E13").Select()
